$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the existing (now-stale) "_GoBack" bookmark that sits in
#    the empty paragraph right after the "options)" syntax line.
#    It lies after the sentence we are about to edit, so removing it
#    first does not disturb earlier character offsets.
# ------------------------------------------------------------------
$old = $d.Bookmarks("_GoBack")
$old.Delete()

# ------------------------------------------------------------------
# 2) Locate the sentence that needs the bookmark wrapped around it
#    (everything except the trailing period).
# ------------------------------------------------------------------
$found = $d.Content
$found.Find.Execute(
    "sets colors of cells in table according to values in the cells relative to the distribution of values.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

# Range covering the sentence without its trailing full stop.
$target = $d.Range($found.Start, $found.End - 1)

# ------------------------------------------------------------------
# 3) Re-add the "_GoBack" bookmark around that sentence (sans the
#    period), which splits the original run into "sentence" +
#    bookmarkEnd + "." just like the tracked edit.
# ------------------------------------------------------------------
$d.Bookmarks.Add("_GoBack", $target)
